$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: G2 "Latest HO Xliff Generate Date"
$wsOverview.Range("G2").Value = "2016-08-31 13:22:19"

# zh-cn sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime"
$wsZhCn.Range("H2").Value = "2016-08-31 13:22:10"
$wsZhCn.Range("K2").Value = "2016-08-31 13:22:41"

# de-de sheet: H2 "Correspond Handoff Datetime", K2 "Correspond Handback DateTime"
$wsDeDe.Range("H2").Value = "2016-08-31 13:22:19"
$wsDeDe.Range("K2").Value = "2016-08-31 13:22:48"
